# Update "Datos actualizados" timestamp and refresh several province figures
# for the provincias_spain.xlsx workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" banner text in A1 (17:05 -> 17:35)
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 17:35"

# Update Casos totales (column B) and Muertes (column E) figures for the
# affected provinces/regions.
$ws.Range("B4").Value = 66643
$ws.Range("E4").Value = 8894

$ws.Range("B5").Value = 55825
$ws.Range("E5").Value = 5981

$ws.Range("B6").Value = 18549
$ws.Range("E6").Value = 1960

$ws.Range("B7").Value = 16677
$ws.Range("E7").Value = 2900

$ws.Range("B9").Value = 12471

$ws.Range("B14").Value = 5520
$ws.Range("E14").Value = 843

$ws.Range("B16").Value = 5202

$ws.Range("B20").Value = 4027
$ws.Range("E20").Value = 352

$ws.Range("B32").Value = 2373
$ws.Range("E32").Value = 303

$ws.Range("E33").Value = 155
